$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the newly-added cells (L6:O6) pick up the same "text" number
# format that the existing barcode columns (B6:K6) already use (style
# index 2 / numFmtId 49) so the long numeric-looking strings are stored
# as text rather than being coerced into floating point numbers.
$ws.Range("L6:O6").NumberFormat = $ws.Range("K6").NumberFormat

# Row 6 ("codigo barras") gets a new set of barcode values; some columns
# repeat the same code (shared string reuse), and four new columns
# (L6:O6) are appended. Set C6 before B6 so the shared-string table ends
# up in the same order as the target file.
$ws.Range("C6").Value = "4157707229253257802023752304390000015702209620171030"
$ws.Range("B6").Value = "4157707229253257802023752304390000005702209620171030"
$ws.Range("D6").Value = "4157707229253257802023752304390000000702209620171030"
$ws.Range("E6").Value = "4157707229253257802023752304390000000552209620171030"
$ws.Range("F6").Value = "4157707229253257802023752304390000015702209620171030"
$ws.Range("G6").Value = "41577072292532578020237523043900000008802209620171030"
$ws.Range("H6").Value = "4157707229253257802023752304390000005712209620171030"
$ws.Range("I6").Value = "4157707229253257802023752304390000015702209620171030"
$ws.Range("J6").Value = "415770722925325780202375230439000007702209620171030"
$ws.Range("K6").Value = "4157707229253257802023752304390000008702209620171030"
$ws.Range("L6").Value = "4157707229253257802023752304390000015702209620171030"
$ws.Range("M6").Value = "4157707229253257802023752304390000015702209620171030"
$ws.Range("N6").Value = "4157707229253257802023752304390000015702209620171030"
$ws.Range("O6").Value = "4157707229253257802023752304390000015702209620171030"

# Selection moved from G10 to H6.
$ws.Range("H6").Select() | Out-Null
